$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.306.95'
$ws.Range("E2").Value = '  +0.26%  '
$ws.Range("D3").Value = '2.606.03'
$ws.Range("E3").Value = '  +0.35%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '537.51'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.33%  '
$ws.Range("E6").Value = '  +0.71%  '
$ws.Range("E7").Value = '  +0.33%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.22%  '
$ws.Range("D9").Value = '2.614.21'
$ws.Range("E9").Value = '  -0.40%  '
$ws.Range("E10").Value = '  -0.64%  '
$ws.Range("E11").Value = '  +1.57%  '
$ws.Range("E12").Value = '  +1.23%  '
$ws.Range("E13").Value = '  +1.83%  '
$ws.Range("D14").Value = '3.067.72'
$ws.Range("E14").Value = '  +0.14%  '
$ws.Range("D15").Value = '59.219.93'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.72%  '
$ws.Range("D17").Value = '2.637.95'
$ws.Range("E17").Value = '  +0.60%  '
$ws.Range("E18").Value = '  +0.65%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '342.39'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.85%  '
$ws.Range("E20").Value = '  +0.75%  '
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("E22").Value = '  -1.58%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.05%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '67.45'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.64%  '
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("E26").Value = '  +0.90%  '
$ws.Range("E27").Value = '  +0.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.52%  '
$ws.Range("E29").Value = '  +0.19%  '
$ws.Range("D30").Value = '0.0₃0738'
$ws.Range("E30").Value = '  +1.72%  '
$ws.Range("E31").Value = '  +5.12%  '
$ws.Range("E32").Value = '  -1.96%  '
$ws.Range("E33").Value = '  -0.13%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.90'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.31%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.96'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.61%  '
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.90'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.49%  '
$ws.Range("E38").Value = '  +0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.825'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.61%  '
$ws.Range("E41").Value = '  +0.24%  '
$ws.Range("E42").Value = '  +0.45%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '274.07'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.61%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.596'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.49%  '
$ws.Range("B45").Value = 'WhiteBITCoin'
$ws.Range("C45").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '10.75'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0961'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.17%  '
$ws.Range("E47").Value = '  +0.54%  '
$ws.Range("D48").Value = '1.944.65'
$ws.Range("E48").Value = '  -2.08%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0222'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.11%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '18.26'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.06%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.49'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.35%  '
